$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column indexes: D=4, J=10, K=11, L=12, M=13, N=14, O=15, P=16, Q=17
$cols = @{ D = 4; J = 10; K = 11; L = 12; M = 13; N = 14; O = 15; P = 16; Q = 17 }

# New values per row for columns D, J, K, L, M, N, O, P, Q
$data = @{
  3  = @{ D = 44186; J = 15; K = 7000;  L = 7000;  M = 7000;  N = "`$/caja 60 unidades"; O = "Provincia de Limarí";   P = 117; Q = 60 }
  4  = @{ D = 44284; J = 35; K = 10000; L = 10000; M = 10000; N = "`$/caja 60 unidades"; O = "Provincia de Limarí";   P = 167; Q = 60 }
  5  = @{ D = 44179; J = 15; K = 7000;  L = 7000;  M = 7000;  N = "`$/caja 60 unidades"; O = "Provincia de Limarí";   P = 117; Q = 60 }
  6  = @{ D = 44291; J = 20; K = 9000;  L = 9000;  M = 9000;  N = "`$/caja 60 unidades"; O = "Provincia de Limarí";   P = 150; Q = 60 }
  7  = @{ D = 44312; J = 30; K = 10000; L = 10000; M = 10000; N = "`$/caja 60 unidades"; O = "Provincia de Limarí";   P = 167; Q = 60 }
  8  = @{ D = 44243; J = 80; K = 10000; L = 11000; M = 10375; N = "`$/caja 60 unidades"; O = "Provincia de Quillota"; P = 173; Q = 60 }
  9  = @{ D = 44315; J = 25; K = 10000; L = 10000; M = 10000; N = "`$/caja 60 unidades"; O = "Provincia de Limarí";   P = 167; Q = 60 }
  11 = @{ D = 44405; J = 45; K = 9000;  L = 9000;  M = 9000;  N = "`$/caja 50 unidades"; O = "Provincia de Quillota"; P = 180; Q = 50 }
  12 = @{ D = 44277; J = 25; K = 10000; L = 10000; M = 10000; N = "`$/caja 60 unidades"; O = "Provincia de Limarí";   P = 167; Q = 60 }
  13 = @{ D = 44585; J = 30; K = 11000; L = 11000; M = 11000; N = "`$/caja 60 unidades"; O = "Provincia de Limarí";   P = 183; Q = 60 }
}

foreach ($row in $data.Keys) {
  $rowData = $data[$row]
  foreach ($colName in $rowData.Keys) {
    $colIndex = $cols[$colName]
    $ws.Cells.Item($row, $colIndex).Value = $rowData[$colName]
  }
}
